$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply cell updates per the commit diff (cryptos list refresh, 2023-12-08).
# D-column values that look like plain numbers get an explicit Text number
# format first so Excel keeps them as strings (preserving things like
# leading/trailing zeros: '0.660', '2.50', etc.) instead of silently
# coercing them to floating point numbers.

$ws.Cells.Item(2, 4).Value = '43.793.73'
$ws.Cells.Item(2, 5).Value = '  +1.14%  '
$ws.Cells.Item(3, 4).Value = '2.368.55'
$ws.Cells.Item(3, 5).Value = '  +5.19%  '
$ws.Cells.Item(4, 5).Value = '  -0.09%  '
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = '0.660'
$ws.Cells.Item(5, 5).Value = '  +2.84%  '
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = '235.33'
$ws.Cells.Item(6, 5).Value = '  +1.91%  '
$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = '73.78'
$ws.Cells.Item(7, 5).Value = '  +13.94%  '
$ws.Cells.Item(8, 5).Value = '  -0.08%  '
$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = '0.522'
$ws.Cells.Item(9, 5).Value = '  +18.65%  '
$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = '0.0989'
$ws.Cells.Item(10, 5).Value = '  +3.79%  '
$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = '27.45'
$ws.Cells.Item(11, 5).Value = '  +2.01%  '
$ws.Cells.Item(12, 4).Value = '2.724.28'
$ws.Cells.Item(12, 5).Value = '  +5.40%  '
$ws.Cells.Item(13, 5).Value = '  +2.48%  '
$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = '16.45'
$ws.Cells.Item(14, 5).Value = '  +9.91%  '
$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = '6.63'
$ws.Cells.Item(15, 5).Value = '  +9.93%  '
$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = '0.886'
$ws.Cells.Item(17, 4).Value = '2.370.16'
$ws.Cells.Item(17, 5).Value = '  +5.33%  '
$ws.Cells.Item(18, 4).Value = '43.681.28'
$ws.Cells.Item(18, 5).Value = '  +1.16%  '
$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = '0.0000101'
$ws.Cells.Item(19, 5).Value = '  +5.43%  '
$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = '75.95'
$ws.Cells.Item(20, 5).Value = '  +3.79%  '
$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = '6.45'
$ws.Cells.Item(21, 5).Value = '  +6.10%  '
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = '250.86'
$ws.Cells.Item(22, 5).Value = '  +1.84%  '
$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = '3.78'
$ws.Cells.Item(24, 5).Value = '  +0.41%  '
$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = '2.50'
$ws.Cells.Item(25, 5).Value = '  +3.47%  '
$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = '10.26'
$ws.Cells.Item(26, 5).Value = '  +5.68%  '
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = '2.25'
$ws.Cells.Item(27, 5).Value = '  -1.61%  '
$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = '22.67'
$ws.Cells.Item(28, 5).Value = '  +4.69%  '
$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = '172.68'
$ws.Cells.Item(29, 5).Value = '  -0.26%  '
$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = '1.54'
$ws.Cells.Item(30, 5).Value = '  +8.42%  '
$ws.Cells.Item(31, 5).Value = '  +3.38%  '
$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = '0.131'
$ws.Cells.Item(32, 5).Value = '  +4.52%  '
$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = '5.13'
$ws.Cells.Item(33, 5).Value = '  +4.20%  '
$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = '0.0703'
$ws.Cells.Item(34, 5).Value = '  +3.81%  '
$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = '5.16'
$ws.Cells.Item(35, 5).Value = '  +5.50%  '
$ws.Cells.Item(36, 5).Value = '  +4.62%  '
$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = '6.67'
$ws.Cells.Item(37, 5).Value = '  +5.20%  '
$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = '2.45'
$ws.Cells.Item(38, 5).Value = '  +7.96%  '
$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = '0.0261'
$ws.Cells.Item(39, 5).Value = '  +4.84%  '
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = '19.67'
$ws.Cells.Item(40, 5).Value = '  +13.44%  '
$ws.Cells.Item(41, 5).Value = '  +0.07%  '
$ws.Cells.Item(42, 5).Value = '  +1.49%  '
$ws.Cells.Item(43, 2).Value = 'Aave'
$ws.Cells.Item(43, 3).Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = '100.67'
$ws.Cells.Item(43, 5).Value = '  +3.98%  '
$ws.Cells.Item(44, 2).Value = 'ARBITRUM'
$ws.Cells.Item(44, 3).Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = '1.17'
$ws.Cells.Item(44, 5).Value = '  +9.66%  '
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = '4.52'
$ws.Cells.Item(45, 5).Value = '  +1.17%  '
$ws.Cells.Item(46, 5).Value = '  +3.40%  '
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = '0.0966'
$ws.Cells.Item(47, 5).Value = '  +3.51%  '
$ws.Cells.Item(48, 4).Value = '1.446.18'
$ws.Cells.Item(48, 5).Value = '  +1.14%  '
$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = '0.176'
$ws.Cells.Item(49, 5).Value = '  +9.75%  '
$ws.Cells.Item(50, 4).Value = '2.592.67'
$ws.Cells.Item(50, 5).Value = '  +4.83%  '
$ws.Cells.Item(51, 2).Value = 'NEARProtocol'
$ws.Cells.Item(51, 3).Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = '2.29'
$ws.Cells.Item(51, 5).Value = '  +0.43%  '
